# The "E_map" column is being split into two erosion-map columns
# (E_map1 / E_map2). Insert a new column at G, duplicating the existing
# E_map (F) values into it, which pushes the former E_exp_Z / E_inv_exp_Z
# columns one slot to the right (G->H, H->I) and extends the used range
# from H180 to I180.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before the old G (E_exp_Z), shifting
# E_exp_Z -> H and E_inv_exp_Z -> I. Column widths/styles of the
# existing columns shift along with it automatically.
$ws.Columns("G").Insert()

# Relabel the headers: F was "E_map", now becomes "E_map1"; the new
# column G becomes "E_map2". H/I keep their shifted-in header text
# (E_exp_Z / E_inv_exp_Z) already, nothing else to do there.
$ws.Cells.Item(1, 6).Value2 = "E_map1"
$ws.Cells.Item(1, 7).Value2 = "E_map2"

# Fill the new E_map2 column with the same values as E_map1 for every
# data row.
for ($r = 2; $r -le 180; $r++) {
    $ws.Cells.Item($r, 7).Value2 = $ws.Cells.Item($r, 6).Value2
}
